$wb = $excel.ActiveWorkbook

# The edit targets the "Subjects_ROC" worksheet (the one with the ROC/NCS
# subject-level table). Grab it explicitly by name so this is correct
# regardless of which sheet happens to be active when the script runs.
$ws = $wb.Worksheets.Item("Subjects_ROC")
$ws.Activate() | Out-Null

# Clear the stale Age (H), BMI (I) and Education (J) values for rows 3-101.
# ClearContents() wipes the value/type but preserves the existing cell
# formatting (the "s" style index), matching the target edit.
$ws.Range("H3:J101").ClearContents() | Out-Null

# Update the worksheet's remembered selection to match the edited range.
$ws.Range("H3:J101").Select() | Out-Null
